# Server Security Setup Input Variables.xlsx - add Perf/QA/PROD server rows
# and a VIM extract security-group user, matching PR 1752.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Server List")
$ws2 = $wb.Worksheets.Item("Security Groups")

# ---------------------------------------------------------------------------
# 1) Security Groups sheet - three new rows for the VIM extract service user
#    (written in the same cell order the original author used, so that new
#    shared strings land in the same sequence as the source workbook).
# ---------------------------------------------------------------------------
$ws2.Cells.Item(41, 4).Value = "Etl"
$ws2.Cells.Item(41, 1).Value = "QA1"
$ws2.Cells.Item(41, 2).Value = "Administrators"
$ws2.Cells.Item(41, 3).Value = "WFM\VIMExtractQA"

$ws2.Cells.Item(42, 1).Value = "Test1"
$ws2.Cells.Item(42, 2).Value = "Administrators"
$ws2.Cells.Item(42, 3).Value = "WFM\VIMExtractTest"
$ws2.Cells.Item(42, 4).Value = "Etl"

$ws2.Cells.Item(43, 1).Value = "Production"
$ws2.Cells.Item(43, 2).Value = "Administrators"
$ws2.Cells.Item(43, 3).Value = "WFM\VIMExtractPrd"
$ws2.Cells.Item(43, 4).Value = "Etl"

# ---------------------------------------------------------------------------
# 2) Server List sheet - add new Perf1 (QA performance) and PROD server rows.
#    Column C gets the same "font1 + border + vertical-center + wrap" look
#    already used lower in the sheet, so copy that format down first.
# ---------------------------------------------------------------------------
$ws1.Cells.Item(24, 3).Copy()
$ws1.Range("C38:C89").PasteSpecial(-4122)
$ws1.Range("C38:C89").WrapText = $true

$webQa = "IRMAQA-PerfWeb07","IRMAQA-PerfWeb08","IRMAQA-PerfWeb09","IRMAQA-PerfWeb10","IRMAQA-PerfWeb11","IRMAQA-PerfWeb12","IRMAQA-PerfWeb01","IRMAQA-PerfWeb02","IRMAQA-PerfWeb03","IRMAQA-PerfWeb04","IRMAQA-PerfWeb05","IRMAQA-PerfWeb06"
for ($i = 0; $i -lt $webQa.Length; $i++) {
    $ws1.Cells.Item(38 + $i, 3).Value = $webQa[$i]
}
$ws1.Cells.Item(38, 1).Value = "Perf1"

$jobQa = "IRMAQA-PerfJob06","IRMAQA-PerfJob07","IRMAQA-PerfJob08","IRMAQA-PerfJob09","IRMAQA-PerfJob10","IRMAQA-PerfJob01","IRMAQA-PerfJob02","IRMAQA-PerfJob03","IRMAQA-PerfJob04","IRMAQA-PerfJob05"
for ($i = 0; $i -lt $jobQa.Length; $i++) {
    $ws1.Cells.Item(50 + $i, 3).Value = $jobQa[$i]
}

$ws1.Cells.Item(60, 3).Value = "IRMAQA-PerfETL02"
$ws1.Cells.Item(61, 3).Value = "IRMAQA-PerfETL01"

$ws1.Cells.Item(62, 2).Value = "TIDAL"
$ws1.Cells.Item(63, 2).Value = "TIDAL"
$ws1.Cells.Item(62, 3).Value = "IRMAQA-PerfTidal02"
$ws1.Cells.Item(63, 3).Value = "IRMAQA-PerfTidal01"

$webPrd = "IRMAPrdWeb01","IRMAPrdWeb07","IRMAPrdWeb02","IRMAPrdWeb08","IRMAPrdWeb03","IRMAPrdWeb09","IRMAPrdWeb04","IRMAPrdWeb10","IRMAPrdWeb05","IRMAPrdWeb11","IRMAPrdWeb06","IRMAPrdWeb12"
for ($i = 0; $i -lt $webPrd.Length; $i++) {
    $ws1.Cells.Item(64 + $i, 3).Value = $webPrd[$i]
}

$jobPrd = "IRMAPrdJob01","IRMAPrdJob06","IRMAPrdJob02","IRMAPrdJob07","IRMAPrdJob03","IRMAPrdJob08","IRMAPrdJob04","IRMAPrdJob09","IRMAPrdJob5","IRMAPrdJob10"
for ($i = 0; $i -lt $jobPrd.Length; $i++) {
    $ws1.Cells.Item(76 + $i, 3).Value = $jobPrd[$i]
}

$ws1.Cells.Item(87, 3).Value = "IRMAPrdETL02"
$ws1.Cells.Item(86, 3).Value = "IRMAPrdETL01"
$ws1.Cells.Item(88, 3).Value = "IRMAPrdTidal01"
$ws1.Cells.Item(89, 3).Value = "IRMAPrdTidal02"

$ws1.Cells.Item(64, 1).Value = "PROD"

# Fill the remaining Env (A) / Server Type (B) columns for the new block.
for ($r = 39; $r -le 63; $r++) { $ws1.Cells.Item($r, 1).Value = "Perf1" }
for ($r = 65; $r -le 89; $r++) { $ws1.Cells.Item($r, 1).Value = "PROD" }

for ($r = 38; $r -le 49; $r++) { $ws1.Cells.Item($r, 2).Value = "Web" }
for ($r = 50; $r -le 59; $r++) { $ws1.Cells.Item($r, 2).Value = "Job" }
for ($r = 60; $r -le 61; $r++) { $ws1.Cells.Item($r, 2).Value = "ETL" }
for ($r = 64; $r -le 75; $r++) { $ws1.Cells.Item($r, 2).Value = "Web" }
for ($r = 76; $r -le 85; $r++) { $ws1.Cells.Item($r, 2).Value = "Job" }
for ($r = 86; $r -le 87; $r++) { $ws1.Cells.Item($r, 2).Value = "ETL" }
$ws1.Cells.Item(88, 2).Value = "TIDAL"
$ws1.Cells.Item(89, 2).Value = "Tidal"

# ---------------------------------------------------------------------------
# 3) Row 36/37 in Server List picked up the column-default border style when
#    the sheet was re-saved; mirror that (style used by C1..C23) onto them.
# ---------------------------------------------------------------------------
$ws1.Cells.Item(2, 3).Copy()
$ws1.Cells.Item(36, 3).PasteSpecial(-4122)
$ws1.Cells.Item(2, 3).Copy()
$ws1.Cells.Item(37, 3).PasteSpecial(-4122)
